$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.677.60"
$ws.Range("E2").Value = "  +1.42%  "

$ws.Range("D3").Value = "1.605.87"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").Value = "'0.518"
$ws.Range("E6").Value = "  +1.53%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").Value = "'27.92"
$ws.Range("E8").Value = "  +5.82%  "

$ws.Range("E9").Value = "  +1.89%  "

$ws.Range("E10").Value = "  +1.74%  "

$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  +0.51%  "

$ws.Range("D12").Value = "1.835.49"
$ws.Range("E12").Value = "  +1.53%  "

$ws.Range("D13").Value = "1.602.98"
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("D14").Value = "'0.549"
$ws.Range("E14").Value = "  +4.83%  "

$ws.Range("D15").Value = "29.688.34"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("E16").Value = "  +0.88%  "

$ws.Range("D17").Value = "'64.06"
$ws.Range("E17").Value = "  +1.86%  "

$ws.Range("D18").Value = "'241.26"
$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("E19").Value = "  +4.61%  "

$ws.Range("D20").Value = "0.0₃0698"
$ws.Range("E20").Value = "  +1.53%  "

$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("E22").Value = "  +0.78%  "

$ws.Range("D23").Value = "'9.40"
$ws.Range("E23").Value = "  +1.79%  "

$ws.Range("E24").Value = "  +0.30%  "

$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("D26").Value = "'15.47"
$ws.Range("E26").Value = "  +2.12%  "

$ws.Range("E27").Value = "  +0.72%  "

$ws.Range("E28").Value = "  +1.69%  "

$ws.Range("E29").Value = "  +0.31%  "

$ws.Range("D30").Value = "'0.0482"
$ws.Range("E30").Value = "  +2.50%  "

$ws.Range("E31").Value = "  +0.33%  "

$ws.Range("E32").Value = "  +0.66%  "

$ws.Range("E33").Value = "  +3.44%  "

$ws.Range("D34").Value = "1.427.66"
$ws.Range("E34").Value = "  +0.30%  "

$ws.Range("D35").Value = "'2.93"
$ws.Range("E35").Value = "  +3.14%  "

$ws.Range("D36").Value = "'1.57"
$ws.Range("E36").Value = "  +4.09%  "

$ws.Range("E37").Value = "  -0.86%  "

$ws.Range("D38").Value = "'2.30"
$ws.Range("E38").Value = "  -0.49%  "

$ws.Range("E39").Value = "  +3.02%  "

$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'0.548"
$ws.Range("E40").Value = "  +3.63%  "

$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D41").Value = "'56.88"
$ws.Range("E41").Value = "  +5.69%  "

$ws.Range("D42").Value = "'0.0500"
$ws.Range("E42").Value = "  +6.47%  "

$ws.Range("D43").Value = "'1.96"
$ws.Range("E43").Value = "  -0.17%  "

$ws.Range("D44").Value = "'0.817"
$ws.Range("E44").Value = "  +2.74%  "

$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.35%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'0.991"
$ws.Range("E46").Value = "  +18.26%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'66.42"
$ws.Range("E47").Value = "  +2.98%  "

$ws.Range("E48").Value = "  -0.11%  "

$ws.Range("D49").Value = "1.745.30"
$ws.Range("E49").Value = "  +1.48%  "

$ws.Range("D50").Value = "'86.72"
$ws.Range("E50").Value = "  +1.10%  "

$ws.Range("D51").Value = "0.0₆0105"
$ws.Range("E51").Value = "  +2.38%  "
